# Auto-generated script applying the cell text changes described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.897.93"
$ws.Range("E2").Value = "  -4.70%  "
$ws.Range("D3").Value = "2.220.03"
$ws.Range("E3").Value = "  -5.69%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'243.84"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("E6").Value = "  -6.45%  "
$ws.Range("D7").Value = "'68.64"
$ws.Range("E7").Value = "  -7.36%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = "  -9.29%  "
$ws.Range("D10").Value = "'0.0953"
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("D11").Value = "'58.10"
$ws.Range("E11").Value = "  -3.74%  "
$ws.Range("D12").Value = "'35.50"
$ws.Range("E12").Value = "  +6.68%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "'6.69"
$ws.Range("E14").Value = "  -7.90%  "
$ws.Range("D15").Value = "2.551.71"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "'14.76"
$ws.Range("E16").Value = "  -9.08%  "
$ws.Range("D17").Value = "'0.840"
$ws.Range("E17").Value = "  -7.40%  "
$ws.Range("D18").Value = "2.219.31"
$ws.Range("E18").Value = "  -5.87%  "
$ws.Range("D19").Value = "41.797.83"
$ws.Range("E19").Value = "  -4.85%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  -7.41%  "
$ws.Range("D21").Value = "'72.39"
$ws.Range("E21").Value = "  -7.41%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  -8.23%  "
$ws.Range("D23").Value = "'234.57"
$ws.Range("E23").Value = "  -6.42%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  +9.51%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -4.92%  "
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").Value = "'9.87"
$ws.Range("E29").Value = "  -5.67%  "
$ws.Range("D30").Value = "'170.45"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").Value = "'20.35"
$ws.Range("E31").Value = "  -8.68%  "
$ws.Range("E32").Value = "  -6.21%  "
$ws.Range("E33").Value = "  -7.54%  "
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").Value = "'5.16"
$ws.Range("E35").Value = "  -4.85%  "
$ws.Range("E36").Value = "  -8.34%  "
$ws.Range("D37").Value = "'3.86"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "'22.59"
$ws.Range("E38").Value = "  +17.22%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0276"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.28"
$ws.Range("E40").Value = "  -5.11%  "
$ws.Range("D41").Value = "'5.82"
$ws.Range("E41").Value = "  -9.47%  "
$ws.Range("D42").Value = "'66.02"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").Value = "'4.93"
$ws.Range("E43").Value = "  -10.21%  "
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("D46").Value = "'0.189"
$ws.Range("E46").Value = "  -5.43%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'4.50"
$ws.Range("E48").Value = "  +5.89%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.18"
$ws.Range("E49").Value = "  -4.36%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").Value = "'10.09"
$ws.Range("E50").Value = "  +6.99%  "
$ws.Range("E51").Value = "  -4.85%  "
